$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the stored credentials (row 2 and row 3, columns A/B)
$ws.Range("A2").Value = "yhernandez"
$ws.Range("B2").Value = "Performance23!"
$ws.Range("A3").Value = "yhernandez"
$ws.Range("B3").Value = "Performance23!"

# Give the username column (A2:A3) a distinct monospace look
$rng = $ws.Range("A2:A3")
$rng.Font.Name = "JetBrains Mono"
$rng.Font.Size = 9.8
$rng.VerticalAlignment = -4108

# Leave the selection where the author left it
$ws.Range("C3").Select()
